# FicheCAF-template.docx
#
# Commit: "déterminer la natureLogement quand elle est null selon la
# natureOperation"
#
# Visible/textual effect inside the template body: the two Jinja guard
# conditions that used to invoke `is_residence()` / `is_foyer()` as method
# calls now read them as plain (property-style) attributes -- the trailing
# "()" after each is removed:
#
#   {% if convention.programme.is_residence() %}  ->  {% if convention.programme.is_residence %}
#   {% if convention.programme.is_foyer() %}       ->  {% if convention.programme.is_foyer %}
#
# (Word's Find/Replace works against the logical run text, so this edits
# cleanly even though "convention.programme.is_residence", "()" and
# " %}" each live in their own <w:r> run.)
#
# The "Détails des logements" table is also (re)applied onto the built-in
# "Table Grid" table style.

$d = $word.ActiveDocument

# 1) convention.programme.is_residence() -> convention.programme.is_residence
$d.Content.Find.Execute("convention.programme.is_residence()", $true, $false, $false, $false, $false, $true, 1, $false, "convention.programme.is_residence", 2) | Out-Null

# 2) convention.programme.is_foyer() -> convention.programme.is_foyer
$d.Content.Find.Execute("convention.programme.is_foyer()", $true, $false, $false, $false, $false, $true, 1, $false, "convention.programme.is_foyer", 2) | Out-Null

# 3) Re-affirm the "Détails des logements" table uses the built-in Table Grid style
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $tbl = $d.Tables($i)
    if ($tbl.Style.NameLocal -eq "Table Grid") {
        $tbl.Style = "Table Grid"
    }
}
